$wb = $excel.ActiveWorkbook

# --- Sheet: Mesh Generation ---
$ms = $wb.Worksheets.Item("Mesh Generation")
$ms.Range("C2").Value = 1.65497067009008
$ms.Range("D2").Value = 0.07108120014891028
$ms.Range("E2").Value = 20.19248209986836
$ms.Range("F2").Value = 2.75426946832548
$ms.Range("C3").Value = 9.06062797771455
$ms.Range("D3").Value = 0.02236719988286495
$ms.Range("E3").Value = 131.3533709999174
$ms.Range("F3").Value = 15.41417822432034

# --- Sheet: Repair ---
$ws = $wb.Worksheets.Item("Repair")

# Insert new rows first (top-to-bottom using post-insert row numbers)
# 1) one new row at position 3 ("Number with NMVs")
$ws.Rows.Item(3).Insert()
# 2) five new rows at position 19 ("Mean/Min/Max/StdDev # of NMVs Before", "Mean # of NMVs After")
$ws.Range("A19:A23").EntireRow.Insert()
# 3) eight net-new rows: replace old 4-row Repair Time block (now at 29-32) with 12 rows
#    Insert 8 blank rows right after the existing 4 rows (so 29-32 stay, 33-40 are fresh)
$ws.Range("A33:A40").EntireRow.Insert()
# 4) one new row at position 42 ("Mean Time per Vertex Repair (s)")
$ws.Range("A42").EntireRow.Insert()

# Now (re)write every label/value for rows 1-51 to their final state
$ws.Range("A1").Value = "Count"
$ws.Range("B1").Value = 1741
$ws.Range("A2").Value = "Number with NMEs"
$ws.Range("B2").Value = 1399
$ws.Range("A3").Value = "Number with NMVs"
$ws.Range("B3").Value = 899
$ws.Range("A4").Value = "Mean # of Added Vertices"
$ws.Range("B4").Value = 142.4537622056289
$ws.Range("A5").Value = "Minimum # of Added Vertices"
$ws.Range("B5").Value = 0
$ws.Range("A6").Value = "Maximum # of Added Vertices"
$ws.Range("B6").Value = 3302
$ws.Range("A7").Value = "Standard Deviation # of Added Vertices"
$ws.Range("B7").Value = 288.0614479167753
$ws.Range("A8").Value = "Mean Added Vertices Percentage"
$ws.Range("B8").Value = 0.005294163506840346
$ws.Range("A9").Value = "Mean # of Added Faces"
$ws.Range("B9").Value = 79.0120620333142
$ws.Range("A10").Value = "Minimum # of Added Faces"
$ws.Range("B10").Value = 0
$ws.Range("A11").Value = "Maximum # of Added Faces"
$ws.Range("B11").Value = 2476
$ws.Range("A12").Value = "Standard Deviation # of Added Faces"
$ws.Range("B12").Value = 188.1529141318901
$ws.Range("A13").Value = "Mean Added Faces Percentage"
$ws.Range("B13").Value = 0.001385859227250776
$ws.Range("A14").Value = "Mean # of NMEs Before"
$ws.Range("B14").Value = 96.4744399770247
$ws.Range("A15").Value = "Minimum # of NMEs Before"
$ws.Range("B15").Value = 0
$ws.Range("A16").Value = "Maximum # of NMEs Before"
$ws.Range("B16").Value = 1929
$ws.Range("A17").Value = "Standard Deviation # of NMEs Before"
$ws.Range("B17").Value = 185.8083788599272
$ws.Range("A18").Value = "Mean # of NMEs After"
$ws.Range("B18").Value = 0
$ws.Range("A19").Value = "Mean # of NMVs Before"
$ws.Range("B19").Value = 20.98506605399196
$ws.Range("A20").Value = "Minimum # of NMVs Before"
$ws.Range("B20").Value = 0
$ws.Range("A21").Value = "Maximum # of NMVs Before"
$ws.Range("B21").Value = 624
$ws.Range("A22").Value = "Standard Deviation # of NMVs Before"
$ws.Range("B22").Value = 47.63407337631442
$ws.Range("A23").Value = "Mean # of NMVs After"
$ws.Range("B23").Value = 0
$ws.Range("A24").Value = "Mean Volume"
$ws.Range("B24").Value = 255851.5502584722
$ws.Range("A25").Value = "Minimum Volume"
$ws.Range("B25").Value = 359
$ws.Range("A26").Value = "Maximum Volume"
$ws.Range("B26").Value = 9533235
$ws.Range("A27").Value = "Standard Deviation Volume"
$ws.Range("B27").Value = 821221.1295896077
$ws.Range("A28").Value = "Mean Volume Change"
$ws.Range("B28").Value = 0
$ws.Range("A29").Value = "Edges Mean Repair Time (s)"
$ws.Range("B29").Value = 0.3515774603990346
$ws.Range("A30").Value = "Edges Minimum Repair Time (s)"
$ws.Range("B30").Value = 0.004164200043305755
$ws.Range("A31").Value = "Edges Maximum Repair Time (s)"
$ws.Range("B31").Value = 9.458849499933422
$ws.Range("A32").Value = "Edges Standard Deviation Repair Time (s)"
$ws.Range("B32").Value = 0.6631316613533466
$ws.Range("A33").Value = "Vertices Mean Repair Time (s)"
$ws.Range("B33").Value = 0.04517836039818542
$ws.Range("A34").Value = "Vertices Minimum Repair Time (s)"
$ws.Range("B34").Value = 0.001613700063899159
$ws.Range("A35").Value = "Vertices Maximum Repair Time (s)"
$ws.Range("B35").Value = 0.6624916999135166
$ws.Range("A36").Value = "Vertices Standard Deviation Repair Time (s)"
$ws.Range("B36").Value = 0.072328306924953
$ws.Range("A37").Value = "Total Mean Repair Time (s)"
$ws.Range("B37").Value = 0.3931456988573239
$ws.Range("A38").Value = "Total Minimum Repair Time (s)"
$ws.Range("B38").Value = 0.006134700030088425
$ws.Range("A39").Value = "Total Maximum Repair Time (s)"
$ws.Range("B39").Value = 9.736861299956217
$ws.Range("A40").Value = "Total Standard Deviation Repair Time (s)"
$ws.Range("B40").Value = 0.7068926608883704
$ws.Range("A41").Value = "Mean Time per Edge Repair (s)"
$ws.Range("B41").Value = 0.002928381819091518
$ws.Range("A42").Value = "Mean Time per Vertex Repair (s)"
$ws.Range("B42").Value = 0.001111683207827253
$ws.Range("A43").Value = "Mean Repair Time Relative to Cuberille Mesh Generation Time"
$ws.Range("B43").Value = 0.2649837108516612
$ws.Range("A44").Value = "Min Repair Time Relative to Cuberille Mesh Generation Time"
$ws.Range("B44").Value = 0.009768309511608573
$ws.Range("A45").Value = "Max Repair Time Relative to Cuberille Mesh Generation Time"
$ws.Range("B45").Value = 2.409305717513604
$ws.Range("A46").Value = "Standard Deviation Repair Time Relative to Cuberille Mesh Generation Time"
$ws.Range("B46").Value = 0.2651064129487007
$ws.Range("A47").Value = "Mean Repair Time Relative to Surface Nets Mesh Generation Time"
$ws.Range("B47").Value = 0.0482680616815102
$ws.Range("A48").Value = "Min Repair Time Relative to Surface Nets Mesh Generation Time"
$ws.Range("B48").Value = 0.001917253643808106
$ws.Range("A49").Value = "Max Repair Time Relative to Surface Nets Mesh Generation Time"
$ws.Range("B49").Value = 0.7558407705526704
$ws.Range("A50").Value = "Standard Deviation Repair Time Relative to Surface Nets Mesh Generation Time"
$ws.Range("B50").Value = 0.07401962938896677
$ws.Range("A51").Value = "Success Rate"
$ws.Range("B51").Value = 1

# Ensure percentage-format rows keep the percent number format / style
$pctCells = @("B8","B13","B43","B44","B45","B46","B47","B48","B49","B50","B51")
foreach ($addr in $pctCells) {
    $ws.Range($addr).NumberFormat = "0.00%"
}